$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-3: account holder name / card number
$ws.Range("C2").Value = "Hartmut"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 21.04.2025"

# Row 6
$ws.Range("B6").Value = "24.04."
$ws.Range("C6").Value = "25.04."
$ws.Range("D6").Value = "AMAZON.DE MKTPLC EU LEENXK"
$ws.Range("E6").Value = "247,65-"

# Row 7
$ws.Range("B7").Value = "25.04."
$ws.Range("C7").Value = "26.04."
$ws.Range("D7").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E7").Value = "25,42-"

# Row 8
$ws.Range("B8").Value = "26.04."
$ws.Range("C8").Value = "27.04."
$ws.Range("D8").Value = "RECHNUNG VODAFONE GMBH 67349"
$ws.Range("E8").Value = "42,23-"

# Row 9: transaction removed -> clear contents, adjust E9 alignment (center/vcenter/wrap)
$ws.Range("B9:D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("E9").HorizontalAlignment = -4108
$ws.Range("E9").VerticalAlignment = -4108
$ws.Range("E9").WrapText = $true

# Row 10: transaction removed -> clear contents, adjust E10 alignment (right/vcenter/wrap)
$ws.Range("B10:D10").ClearContents()
$ws.Range("E10").ClearContents()
$ws.Range("E10").HorizontalAlignment = -4152
$ws.Range("E10").VerticalAlignment = -4108
$ws.Range("E10").WrapText = $true

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 28.04.2025"
$ws.Range("E12").Value = "315,30-"

# Next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 07.05.2025"
